# Applies the "Update metadata to reflect the new state" change to the
# "events" sheet: two table columns are renamed, a new regular/extension
# visit row is introduced, and a couple of workbook-level selections move.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "events" worksheet: update the event-pattern table
# ---------------------------------------------------------------------------
$wsEvents = $wb.Worksheets.Item("events")
$lo = $wsEvents.ListObjects.Item(1)

# Grow the table by one row (A1:F7 -> A1:F8) and give the table its new name.
$lo.ListRows.Add() | Out-Null
$lo.Name = "Table91310"

# Rename the two table columns that changed meaning (renaming the header
# cell itself is what actually updates the table's column definition).
$wsEvents.Range("C1").Value = "is_regular_visit"
$wsEvents.Range("F1").Value = "is_baseline_event"

# Clear out all the old data rows so they can be rewritten cleanly.
$wsEvents.Range("A2:F8").ClearContents()

# Re-populate the table with the updated event definitions.
# event_id | event_id_pattern | is_regular_visit | event_label_custom | event_name_custom | is_baseline_event
$wsEvents.Range("A2").Value = "SCR"
$wsEvents.Range("C2").Value = $true
$wsEvents.Range("E2").Value = "Screening"
$wsEvents.Range("F2").Value = $true

$wsEvents.Range("B3").Value = "^VIS$"
$wsEvents.Range("C3").Value = $true
$wsEvents.Range("D3").Value = "V"
$wsEvents.Range("E3").Value = "Visit"
$wsEvents.Range("F3").Value = $false

$wsEvents.Range("B4").Value = "^VISEXT"
$wsEvents.Range("C4").Value = $false
$wsEvents.Range("E4").Value = "Ext. visit"
$wsEvents.Range("F4").Value = $false

$wsEvents.Range("A5").Value = "EOT"
$wsEvents.Range("C5").Value = $true
$wsEvents.Range("D5").Value = "EoT"
$wsEvents.Range("E5").Value = "EoT"
$wsEvents.Range("F5").Value = $false

$wsEvents.Range("A6").Value = "FU1"
$wsEvents.Range("C6").Value = $true
$wsEvents.Range("D6").Value = "FU"
$wsEvents.Range("E6").Value = "Visit"
$wsEvents.Range("F6").Value = $false

$wsEvents.Range("A7").Value = "EXIT"
$wsEvents.Range("C7").Value = $false
$wsEvents.Range("E7").Value = "Exit"
$wsEvents.Range("F7").Value = $false

$wsEvents.Range("B8").Value = "^UN"
$wsEvents.Range("C8").Value = $false
$wsEvents.Range("E8").Value = "Unscheduled visit"
$wsEvents.Range("F8").Value = $false

# ---------------------------------------------------------------------------
# 2. Update the selection shown on "study_forms" (without changing which
#    sheet is active).
# ---------------------------------------------------------------------------
$wsStudyForms = $wb.Worksheets.Item("study_forms")
$wsStudyForms.Activate()
$wsStudyForms.Range("A10").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Make "events" the active sheet/tab, with D3 selected.
# ---------------------------------------------------------------------------
$wsEvents.Activate()
$wsEvents.Range("D3").Select() | Out-Null
